$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the shared-string text values used across the "Results"/"Remarks"
# columns. The workbook previously had:
#   FAIL
#   SCANNER DID NOT CATCH THE ERROR
#   line 12 input accepted
# and now the text has been edited in-place (re-using cells on rows 15/23)
# per Alla's note: the bug was fixed on 4/23 by Jared Cox.

# Row 15 (PREREQ013): Results -> "FAIL(PASS NOW)", Remarks -> updated SCANNER text
$ws.Range("D15").Value = "FAIL(PASS NOW)"
$ws.Range("E15").Value = "SCANNER DID NOT CATCH THE ERROR  (FIXED ON 4/23 BY JARED COX)"

# Row 23 (PREREQ021): Results -> "FAIL(PASS NOW)", Remarks -> updated line-12 text
$ws.Range("D23").Value = "FAIL(PASS NOW)"
$ws.Range("E23").Value = "line 12 input accepted (FIXED ON 4/23 BY JARED COX)"

# Both rows now wrap + center their Results/Remarks text (matching the style
# already used on column A of those same rows) and grow to a 30pt row height
# so the longer text fits.
$rng15 = $ws.Range("D15:E15")
$rng15.HorizontalAlignment = -4108
$rng15.VerticalAlignment = -4108
$rng15.WrapText = $true

$rng23 = $ws.Range("D23:E23")
$rng23.HorizontalAlignment = -4108
$rng23.VerticalAlignment = -4108
$rng23.WrapText = $true

$ws.Rows.Item(15).RowHeight = 30
$ws.Rows.Item(23).RowHeight = 30

# Restore the view/selection to match the saved state (scrolled so row 7 is
# at the top, with E23 as the active cell/selection).
$win = $excel.ActiveWindow
$win.ScrollRow = 7
$win.ScrollColumn = 1
$ws.Range("E23").Select()
